# Swap the observation data between row 34 and row 35 for the
# columns that differ between the two records (A, B, E, F, G, H, Q, R, X, Z, AB, AC).
# Columns D, I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 34
$row2 = 35

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "X", "Z", "AB", "AC")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"

    $val1 = $ws.Range($addr1).Value2
    $val2 = $ws.Range($addr2).Value2

    $ws.Range($addr1).Value = $val2
    $ws.Range($addr2).Value = $val1
}
